$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I (dS0 analog) and J (dSF analog)
$iValues = @(7,7,2,6,8,7,7,3,8,6,7,7,9,4,6,8,3,7,3,1,5)
$jValues = @(7,8,2,6,8,8,8,3,8,7,8,9,9,5,7,8,5,7,3,1,5)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
